$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new Table4 row (row 12, columns M:Q) that was blank before.
$ws.Range("M12").Value = "Problem Solving(Algorithms & Data Structures)"
$ws.Range("N12").Value = 45120
$ws.Range("O12").Value = "1172.97/2200"
$ws.Range("P12").Value = 124134
$ws.Range("Q12").Formula = "=IF(ROW()>2,(`$P`$2-P12)/`$P`$2,""NA"")"

# Update the selected cell shown in the sheet view.
$ws.Range("M12").Select()

# Shift the saved window position (simulates the window being moved/resized).
$excel.ActiveWindow.Left = -28920
